$d = $word.ActiveDocument

# Insert " again" right after "Something is changed" text, within the same run/paragraph.
$rng = $d.Content
$rng.Find.Execute("Something is changed", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)  # wdCollapseEnd
$rng.Text = " again"
$rng.Font.Reset()

# Remove the existing _GoBack bookmark from the second paragraph, then recreate it
# at the end of the first paragraph (after the newly inserted text).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$p1 = $d.Paragraphs(1).Range
$p1.Collapse(0)
$p1.MoveEnd(1, -1) | Out-Null
$d.Bookmarks.Add("_GoBack", $p1) | Out-Null
